$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.739
$ws.Range("AF5").Value = 0.967
$ws.Range("AF6").Value = 0.838
$ws.Range("AF7").Value = 0.911
$ws.Range("AF8").Value = 0.876
$ws.Range("AF9").Value = 0.733
$ws.Range("AF10").Value = 0.967
$ws.Range("AF11").Value = 0.967
$ws.Range("AF12").Value = 1.276
$ws.Range("AF13").Value = 1.6
